$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 326
$ws.Range("I5").Value = 326
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 326
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -211
$ws.Range("N5").Value = ""
$ws.Range("H34").Value = 24949.5
$ws.Range("I34").Value = 9900
$ws.Range("K34").Value = 9900
$ws.Range("M34").Value = -9697
$ws.Range("H36").Value = 24949.5
$ws.Range("I36").Value = 9900
$ws.Range("K36").Value = 9900
$ws.Range("M36").Value = -9185
$ws.Range("H62").Value = 6109.8423
$ws.Range("I62").Value = 5463.7144
$ws.Range("K62").Value = 5463.7144
$ws.Range("M62").Value = -4839.7144
$ws.Range("H65").Value = 6109.8423
$ws.Range("I65").Value = 5463.7144
$ws.Range("K65").Value = 27318.572
$ws.Range("M65").Value = -24198.572
$ws.Range("H70").Value = 203699.2
$ws.Range("I70").Value = 4999
$ws.Range("J70").Value = 501749.5
$ws.Range("K70").Value = 14997
$ws.Range("L70").Value = 1505248.5
$ws.Range("M70").Value = -14727
$ws.Range("N70").Value = -1505788.5
$ws.Range("H73").Value = 203699.2
$ws.Range("I73").Value = 4999
$ws.Range("J73").Value = 501749.5
$ws.Range("K73").Value = 14997
$ws.Range("L73").Value = 1505248.5
$ws.Range("M73").Value = -14061
$ws.Range("N73").Value = -1507120.5
$ws.Range("H106").Value = 8074.9287
$ws.Range("I106").Value = 8074.9287
$ws.Range("K106").Value = 8074.9287
$ws.Range("M106").Value = -7443.9287
$ws.Range("H138").Value = 2599.1758
$ws.Range("I138").Value = 2197.6
$ws.Range("J138").Value = 2661.9219
$ws.Range("K138").Value = 6592.799999999999
$ws.Range("L138").Value = 7985.7657
$ws.Range("M138").Value = -1452.799999999999
$ws.Range("N138").Value = -18265.7657

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5448.8
$ws.Range("I20").Value = 5312.375
$ws.Range("J20").Value = 5994.5
$ws.Range("K20").Value = 5312.375
$ws.Range("L20").Value = 5994.5
$ws.Range("M20").Value = -5065.375
$ws.Range("N20").Value = -6488.5
$ws.Range("H86").Value = 6412385
$ws.Range("I86").Value = 11905787
$ws.Range("K86").Value = 11905787
$ws.Range("M86").Value = -11904664
$ws.Range("H89").Value = 6412385
$ws.Range("I89").Value = 11905787
$ws.Range("K89").Value = 59528935
$ws.Range("M89").Value = -59523319
$ws.Range("H94").Value = 31271930
$ws.Range("I94").Value = 39474572
$ws.Range("K94").Value = 39474572
$ws.Range("M94").Value = -39474121
$ws.Range("H99").Value = 1851.1818
$ws.Range("I99").Value = 1851.1818
$ws.Range("K99").Value = 1851.1818
$ws.Range("M99").Value = -353.1818000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 873.5238000000001
$ws.Range("I22").Value = 875.4211
$ws.Range("J22").Value = 855.5
$ws.Range("K22").Value = 875.4211
$ws.Range("L22").Value = 855.5
$ws.Range("M22").Value = -525.4211
$ws.Range("N22").Value = -1555.5
$ws.Range("H62").Value = 7074.0713
$ws.Range("I62").Value = 7096.3
$ws.Range("K62").Value = 7096.3
$ws.Range("M62").Value = -6472.3
$ws.Range("H65").Value = 7074.0713
$ws.Range("I65").Value = 7096.3
$ws.Range("K65").Value = 35481.5
$ws.Range("M65").Value = -32361.5
$ws.Range("H122").Value = 84057.73
$ws.Range("I122").Value = 110108.43
$ws.Range("K122").Value = 330325.29
$ws.Range("M122").Value = -327875.29
$ws.Range("H132").Value = 1641.1
$ws.Range("I132").Value = 1362.9333
$ws.Range("K132").Value = 4088.7999
$ws.Range("M132").Value = -1558.7999
$ws.Range("H134").Value = 1393
$ws.Range("I134").Value = 1397.5
$ws.Range("K134").Value = 4192.5
$ws.Range("M134").Value = -1657.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 212536000
$ws.Range("J4").Value = 20663334
$ws.Range("L4").Value = 61990002
$ws.Range("N4").Value = -61990226
$ws.Range("H7").Value = 7008.6924
$ws.Range("I7").Value = 7880.5557
$ws.Range("K7").Value = 23641.6671
$ws.Range("M7").Value = -23529.6671
$ws.Range("H38").Value = 67.2
$ws.Range("I38").Value = 67.2
$ws.Range("K38").Value = 201.6
$ws.Range("M38").Value = 145.4
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H92").Value = 717.36365
$ws.Range("J92").Value = 717.36365
$ws.Range("L92").Value = 2152.09095
$ws.Range("N92").Value = -4648.09095
$ws.Range("H109").Value = 5053.579
$ws.Range("J109").Value = 5557
$ws.Range("L109").Value = 16671
$ws.Range("N109").Value = -18751
$ws.Range("H131").Value = 2418.85
$ws.Range("I131").Value = 1478.125
$ws.Range("J131").Value = 3046
$ws.Range("K131").Value = 4434.375
$ws.Range("L131").Value = 9138
$ws.Range("M131").Value = 605.625
$ws.Range("N131").Value = -19218
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4249.75
$ws.Range("I70").Value = 2333
$ws.Range("K70").Value = 2333
$ws.Range("M70").Value = -2063
$ws.Range("H73").Value = 4249.75
$ws.Range("I73").Value = 2333
$ws.Range("K73").Value = 2333
$ws.Range("M73").Value = -1397
$ws.Range("H80").Value = 7095.048
$ws.Range("I80").Value = 6888.875
$ws.Range("J80").Value = 7221.923
$ws.Range("K80").Value = 6888.875
$ws.Range("L80").Value = 7221.923
$ws.Range("M80").Value = -5890.875
$ws.Range("N80").Value = -9217.922999999999
$ws.Range("H83").Value = 7095.048
$ws.Range("I83").Value = 6888.875
$ws.Range("J83").Value = 7221.923
$ws.Range("K83").Value = 34444.375
$ws.Range("L83").Value = 36109.615
$ws.Range("M83").Value = -29452.375
$ws.Range("N83").Value = -46093.615
$ws.Range("H107").Value = 750.63635
$ws.Range("I107").Value = 665.7
$ws.Range("K107").Value = 665.7
$ws.Range("M107").Value = 1254.3
$ws.Range("H126").Value = 6456.355
$ws.Range("I126").Value = 5398.0454
$ws.Range("K126").Value = 16194.1362
$ws.Range("M126").Value = -13724.1362
$ws.Range("H136").Value = 30706.625
$ws.Range("J136").Value = 30706.625
$ws.Range("L136").Value = 92119.875
$ws.Range("N136").Value = -97219.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 375
$ws.Range("J22").Value = 450
$ws.Range("L22").Value = 450
$ws.Range("N22").Value = -1040
$ws.Range("H27").Value = 375
$ws.Range("J27").Value = 450
$ws.Range("L27").Value = 450
$ws.Range("N27").Value = -664
$ws.Range("H55").Value = 396
$ws.Range("I55").Value = 396
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 396
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -223
$ws.Range("N55").Value = ""
$ws.Range("H82").Value = 1114.625
$ws.Range("I82").Value = 989.4286
$ws.Range("J82").Value = 1991
$ws.Range("K82").Value = 989.4286
$ws.Range("L82").Value = 1991
$ws.Range("M82").Value = -628.4286
$ws.Range("N82").Value = -2713
$ws.Range("H85").Value = 1114.625
$ws.Range("I85").Value = 989.4286
$ws.Range("J85").Value = 1991
$ws.Range("K85").Value = 989.4286
$ws.Range("L85").Value = 1991
$ws.Range("M85").Value = 258.5714
$ws.Range("N85").Value = -4487
$ws.Range("H93").Value = 1365.68
$ws.Range("I93").Value = 1277.8667
$ws.Range("J93").Value = 1497.4
$ws.Range("K93").Value = 1277.8667
$ws.Range("L93").Value = 1497.4
$ws.Range("M93").Value = -29.86670000000004
$ws.Range("N93").Value = -3993.4
$ws.Range("H132").Value = 2369.3125
$ws.Range("I132").Value = 2369.3125
$ws.Range("K132").Value = 7107.9375
$ws.Range("M132").Value = -4577.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25035
$ws.Range("I54").Value = 25035
$ws.Range("K54").Value = 25035
$ws.Range("M54").Value = -24515
$ws.Range("H62").Value = 8323
$ws.Range("I62").Value = 8387.6
$ws.Range("K62").Value = 8387.6
$ws.Range("M62").Value = -7763.6
$ws.Range("H65").Value = 8323
$ws.Range("I65").Value = 8387.6
$ws.Range("K65").Value = 41938
$ws.Range("M65").Value = -38818
$ws.Range("H126").Value = 7474.75
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 8299.666999999999
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 24899.001
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -29839.001
$ws.Range("H136").Value = 4724.7646
$ws.Range("I136").Value = 4795.5
$ws.Range("J136").Value = 4494.875
$ws.Range("K136").Value = 14386.5
$ws.Range("L136").Value = 13484.625
$ws.Range("M136").Value = -11836.5
$ws.Range("N136").Value = -18584.625

